$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text runs within shared strings) ---
# A8: "Volume 32   Number  23" -> "...24" (chars 21-22, 1-based)
$ws.Cells.Item(8,1).Characters(21,2).Text = "24"

# C9: "Report Covering the Week  6/2/2025  Through  6/8/2025"
#     -> "...6/9/2025  Through  6/15/2025"
$ws.Cells.Item(9,3).Characters(27,8).Text = "6/9/2025"
$ws.Cells.Item(9,3).Characters(46,8).Text = "6/15/2025"

# --- C28: numeric 2 -> text "0" (shared-string, matches D28 exactly) ---
$ws.Range("D28").Copy($ws.Range("C28"))

# --- Crime-stat table numeric updates (rows 15-28) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("L15").Value = -22.222222222222
$ws.Range("N15").Value = -68.181818181818
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 54
$ws.Range("J16").Value = 57
$ws.Range("K16").Value = -5.263157894736
$ws.Range("L16").Value = -32.5
$ws.Range("M16").Value = -22.857142857142
$ws.Range("N16").Value = -85.405405405405
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -41.176470588235
$ws.Range("I17").Value = 86
$ws.Range("J17").Value = 83
$ws.Range("K17").Value = 3.614457831325
$ws.Range("L17").Value = -23.893805309734
$ws.Range("M17").Value = 22.857142857142
$ws.Range("N17").Value = -66.007905138339
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -26.315789473684
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = -1.25
$ws.Range("L18").Value = -38.759689922480
$ws.Range("M18").Value = -24.038461538461
$ws.Range("N18").Value = -77.808988764044
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -23.076923076923
$ws.Range("F19").Value = 62
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = 3.333333333333
$ws.Range("I19").Value = 369
$ws.Range("J19").Value = 333
$ws.Range("K19").Value = 10.810810810810
$ws.Range("L19").Value = -19.432314410480
$ws.Range("M19").Value = 8.529411764705
$ws.Range("N19").Value = -44.925373134328
$ws.Range("D20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -80
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -5.882352941176
$ws.Range("L20").Value = -5.882352941176
$ws.Range("M20").Value = -23.809523809523
$ws.Range("N20").Value = -93.073593073593
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 102
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = -8.108108108108
$ws.Range("I21").Value = 612
$ws.Range("J21").Value = 579
$ws.Range("K21").Value = 5.699481865284
$ws.Range("L21").Value = -24.163568773234
$ws.Range("M21").Value = 0.163666121112
$ws.Range("N21").Value = -67.90770844258
$ws.Range("C23").Value = 6
$ws.Range("E23").Value = 500
$ws.Range("F23").Value = 13
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 63
$ws.Range("J23").Value = 35
$ws.Range("K23").Value = 80
$ws.Range("L23").Value = 1.612903225806
$ws.Range("M23").Value = 21.153846153846
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -41.025641025641
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 135
$ws.Range("H24").Value = -34.814814814814
$ws.Range("I24").Value = 680
$ws.Range("J24").Value = 671
$ws.Range("K24").Value = 1.341281669150
$ws.Range("L24").Value = 4.615384615384
$ws.Range("M24").Value = -6.464924346629
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -47.619047619047
$ws.Range("F25").Value = 40
$ws.Range("G25").Value = 81
$ws.Range("H25").Value = -50.617283950617
$ws.Range("I25").Value = 331
$ws.Range("J25").Value = 389
$ws.Range("K25").Value = -14.910025706940
$ws.Range("L25").Value = -12.894736842105
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 35
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 193
$ws.Range("J26").Value = 173
$ws.Range("K26").Value = 11.560693641618
$ws.Range("L26").Value = -11.059907834101
$ws.Range("M26").Value = -6.763285024154
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = -52.941176470588
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
